# Update the sample data across all three sheets: replace the people
# (Dorka Zomok / Alejandro Novoa / Ryan Boggio) with famous space scientists
# (Carl Sagan / Wernher von Braun / Szergej Koroljov) and fix up their ages.
$wb = $excel.ActiveWorkbook

# --- Sheet1 (rows in natural order: Number 1,2,3) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B2").Value = "Carl"
$ws1.Range("C2").Value = "Sagan"
$ws1.Range("D2").Value = 67
$ws1.Range("B3").Value = "Wernher"
$ws1.Range("C3").Value = "von Braun"
$ws1.Range("D3").Value = 88
$ws1.Range("B4").Value = "Szergej"
$ws1.Range("C4").Value = "Koroljov"
$ws1.Range("D4").Value = 71
$ws1.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws1.Columns.Item(4).ColumnWidth = 3.6666666666666665

# --- Sheet2 (rows in natural order: Number 1,2,3) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B2").Value = "Carl"
$ws2.Range("C2").Value = "Sagan"
$ws2.Range("D2").Value = 67
$ws2.Range("B3").Value = "Wernher"
$ws2.Range("C3").Value = "von Braun"
$ws2.Range("D3").Value = 88
$ws2.Range("B4").Value = "Szergej"
$ws2.Range("C4").Value = "Koroljov"
$ws2.Range("D4").Value = 71
$ws2.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws2.Columns.Item(4).ColumnWidth = 3.6666666666666665

# --- Sheet3 (rows are reordered: Number 2, 3, 1) ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B2").Value = "Wernher"
$ws3.Range("C2").Value = "von Braun"
$ws3.Range("D2").Value = 88
$ws3.Range("B3").Value = "Szergej"
$ws3.Range("C3").Value = "Koroljov"
$ws3.Range("D3").Value = 71
$ws3.Range("B4").Value = "Carl"
$ws3.Range("C4").Value = "Sagan"
$ws3.Range("D4").Value = 67
$ws3.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws3.Columns.Item(4).ColumnWidth = 3.6666666666666665
